$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells in rows 2-19 (columns B-F) per diff
$ws.Range("B2").Value = "NSE:AARON"
$ws.Range("C2").Value = "NSE:ABB"
$ws.Range("D2").Value = "NSE:JINDALSTEL"
$ws.Range("E2").Value = "NSE:PHOENIXLTD"
$ws.Range("F2").Value = "NSE:FORTIS"
$ws.Range("B3").Value = "NSE:DATAMATICS"
$ws.Range("C3").Value = "NSE:AGARIND"
$ws.Range("B4").Value = "NSE:FORTIS"
$ws.Range("C4").Value = "NSE:ANANTRAJ"
$ws.Range("B5").Value = "NSE:GALAXYSURF"
$ws.Range("C5").Value = "NSE:ATL"
$ws.Range("B6").Value = "NSE:GSFC"
$ws.Range("C6").Value = "NSE:ATUL"
$ws.Range("B7").Value = "NSE:GUFICBIO"
$ws.Range("C7").Value = "NSE:BCLIND"
$ws.Range("B8").Value = "NSE:HDFCLIQUID"
$ws.Range("C8").Value = "NSE:CAPACITE"
$ws.Range("B9").Value = "NSE:HLEGLAS"
$ws.Range("C9").Value = "NSE:CDSL"
$ws.Range("B10").Value = "NSE:KAMATHOTEL"
$ws.Range("C10").Value = "NSE:CONTROLPR"
$ws.Range("B11").Value = "NSE:KIMS"
$ws.Range("C11").Value = "NSE:DEN"
$ws.Range("B12").Value = "NSE:LALPATHLAB"
$ws.Range("C12").Value = "NSE:DHANBANK"
$ws.Range("B13").Value = "NSE:NAVINFLUOR"
$ws.Range("C13").Value = "NSE:DHANUKA"
$ws.Range("B14").Value = "NSE:OLECTRA"
$ws.Range("C14").Value = "NSE:DISHTV"
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "NSE:GOLDIAM"
$ws.Range("C16").Value = "NSE:HERITGFOOD"
$ws.Range("C17").Value = "NSE:HPAL"
$ws.Range("C18").Value = "NSE:INDHOTEL"
$ws.Range("C19").Value = "NSE:INDIAGLYCO"

# Copy formatting (style) from A15 (bold, centered, bordered) down to new rows A20:A31
$ws.Range("A15").Copy()
$ws.Range("A20:A31").PasteSpecial(-4122)

# Populate new rows 20-31 (column A index, column C ticker)
$ws.Range("A20").Value = 18
$ws.Range("C20").Value = "NSE:KIRLOSIND"
$ws.Range("A21").Value = 19
$ws.Range("C21").Value = "NSE:KIRLPNU"
$ws.Range("A22").Value = 20
$ws.Range("C22").Value = "NSE:KOKUYOCMLN"
$ws.Range("A23").Value = 21
$ws.Range("C23").Value = "NSE:MASTEK"
$ws.Range("A24").Value = 22
$ws.Range("C24").Value = "NSE:MMTC"
$ws.Range("A25").Value = 23
$ws.Range("C25").Value = "NSE:NARMADA"
$ws.Range("A26").Value = 24
$ws.Range("C26").Value = "NSE:NBCC"
$ws.Range("A27").Value = 25
$ws.Range("C27").Value = "NSE:NGLFINE"
$ws.Range("A28").Value = 26
$ws.Range("C28").Value = "NSE:NILASPACES"
$ws.Range("A29").Value = 27
$ws.Range("C29").Value = "NSE:PPLPHARMA"
$ws.Range("A30").Value = 28
$ws.Range("C30").Value = "NSE:RATEGAIN"
$ws.Range("A31").Value = 29
$ws.Range("C31").Value = "NSE:ROHLTD"

Write-Output "done"
